# payweek / nonpayweek is now displayed on each day of the roster
#
# - Staff sheet header row (D1:W1) is relabelled from "<Day> 1/2 start/end"
#   (week-number based) to "<Day> start/end P/N" (payweek / non-payweek
#   based).
# - The active sheet/selection moves from Shifts!F7 to Staff!H8.

$wb = $excel.ActiveWorkbook

$shifts = $wb.Worksheets.Item("Shifts")
$staff  = $wb.Worksheets.Item("Staff")

# Relabel the fortnightly roster headers on the Staff sheet: the "1"/"2"
# week-number suffixes become "P" (payweek) / "N" (non-payweek) suffixes,
# and "start"/"end" moves in front of the day name.
$staff.Range("D1").Value = "Mon start P"
$staff.Range("E1").Value = "Mon end P"
$staff.Range("F1").Value = "Tue start P"
$staff.Range("G1").Value = "Tue end P"
$staff.Range("H1").Value = "Wed start P"
$staff.Range("I1").Value = "Wed end P"
$staff.Range("J1").Value = "Thu start P"
$staff.Range("K1").Value = "Thu end P"
$staff.Range("L1").Value = "Fri start P"
$staff.Range("M1").Value = "Fri end P"
$staff.Range("N1").Value = "Mon start N"
$staff.Range("O1").Value = "Mon end N"
$staff.Range("P1").Value = "Tue start N"
$staff.Range("Q1").Value = "Tue end N"
$staff.Range("R1").Value = "Wed start N"
$staff.Range("S1").Value = "Wed End N"
$staff.Range("T1").Value = "Thu start N"
$staff.Range("U1").Value = "Thu end N"
$staff.Range("V1").Value = "Fri start N"
$staff.Range("W1").Value = "Fri end N"

# Move the selection/active sheet from Shifts!F7 to Staff!H8.
[void]$shifts.Range("E15").Select()
[void]$staff.Activate()
[void]$staff.Range("H8").Select()
